$wb = $excel.ActiveWorkbook

# --- OFF sheet: Week 15 logged ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 250
$wsOff.Range("C2").Value = 178
$wsOff.Range("D2").Value = 52
$wsOff.Range("G2").Value = 2

# --- DEF sheet: Week 16 simulated ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 167
$wsDef.Range("C2").Value = 111
$wsDef.Range("D2").Value = 43
$wsDef.Range("F2").Value = 4
